# Weekly price-sheet update: a new week of data is inserted at row 55,
# pushing every existing record (old rows 55-125) down by one row
# (new rows 56-126). The new row 55 carries a fresh reading for
# "Feria Lagunitas de Puerto Montt" / "Pepino ensalada".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 55; Excel shifts rows 55:125 down to 56:126 and
# copies formatting (incl. the date style on column D) from the row above.
$ws.Rows.Item(55).Insert()

$ws.Range("A55").Value2 = 4
$ws.Range("B55").Value2 = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C55").Value2 = 'Los Lagos'
$ws.Range("D55").Value2 = 44413
$ws.Range("E55").Value2 = 10
$ws.Range("F55").Value2 = 100112043
$ws.Range("G55").Value2 = 'Pepino ensalada'
$ws.Range("H55").Value2 = 'Sin especificar'
$ws.Range("I55").Value2 = 'Primera'
$ws.Range("J55").Value2 = 250
$ws.Range("K55").Value2 = 19000
$ws.Range("L55").Value2 = 19000
$ws.Range("M55").Value2 = 19000
$ws.Range("N55").Value2 = '$/caja 60 unidades'
$ws.Range("O55").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P55").Value2 = 317
$ws.Range("Q55").Value2 = 60
$ws.Range("R55").Value2 = 'Hortaliza'
